$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1501
$ws1.Range("F5").Value = 206
$ws1.Range("F7").Value = 94
$ws1.Range("F8").Value = 9867
$ws1.Range("F10").Value = 118
$ws1.Range("F13").Value = 374
$ws1.Range("F14").Value = 6827
$ws1.Range("F15").Value = 1084
$ws1.Range("F16").Value = 632
$ws1.Range("F17").Value = 51
$ws1.Range("F18").Value = 197

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1501
$ws4.Range("F5").Value = 206
$ws4.Range("F8").Value = 95
$ws4.Range("F11").Value = 9867
$ws4.Range("F13").Value = 118
$ws4.Range("F16").Value = 374
$ws4.Range("F17").Value = 6827
$ws4.Range("F18").Value = 1084
$ws4.Range("F19").Value = 632
$ws4.Range("F20").Value = 51
$ws4.Range("F21").Value = 197
